$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("by logger")
$ws2 = $wb.Worksheets.Item("by site")

$ws1.Range("A34").Value = "06-01636"
$ws1.Range("A34").Font.Name = "Segoe UI"
$ws1.Range("A34").Font.Size = 11
$ws1.Range("A34").Font.Bold = $false
$ws1.Range("A34").HorizontalAlignment = -4131
$ws1.Range("A34").VerticalAlignment = -4160
$ws1.Range("A34").WrapText = $true

$ws1.Range("B34").Value = "LN-K01-12177-1"
$ws1.Range("C34").Value = "1685401568; 1286503872"

$ws2.Range("A38").Value = "LN-K01-12177-1"
$ws2.Range("B38").Value = "06-01636"
$ws2.Range("C38").Value = 3
$ws2.Range("D38").Value = "5570f23feeb666615003051140cea73ccdb18639"

$ws1.Activate()
